$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new hab entries (appended below the existing data; a
# descending sort on Population further below moves them into their
# correct sorted position). ---

# New hab: Gavaswadi (Orbit / Solaris / Haunting Mars / ORB 2 / Torus)
$ws.Range("A43").Value = "Gavaswadi"
$ws.Range("B43").Value = "Orbit"
$ws.Range("C43").Value = "Solaris"
$ws.Range("D43").Value = "Haunting Mars"
$ws.Range("E43").Value = "ORB 2"
$ws.Range("F43").Value = "Torus"
$ws.Range("G43").Value = 3000
$ws.Range("H43").Value = 0.1
$ws.Range("I43").Formula = "=G43*H43"
$ws.Range("J43").Value = "Orbital infrastructure ops"

# New hab: Piros Lyuk (S/ / Movement / Haunting Mars / Crater)
$ws.Range("A44").Value = "Piros Lyuk"
$ws.Range("B44").Value = "S/"
$ws.Range("C44").Value = "Movement"
$ws.Range("D44").Value = "Haunting Mars"
$ws.Range("F44").Value = "Crater"
$ws.Range("G44").Value = 600
$ws.Range("H44").Value = 0.6
$ws.Range("I44").Formula = "=G44*H44"
$ws.Range("J44").Value = "Self-sufficient terraforming camp"

# --- Re-sort the whole table (A1:N44) descending by Population (col G),
# same as the sheet's existing sort state, so the two new rows land in
# the correct position among the existing, already-sorted rows. ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("G1:G44"), 0, 2, 0, 0)
$sortObj.SetRange($ws.Range("A1:N44"))
$sortObj.Header = 2
$sortObj.Apply()

# The sort shifts the "total" row (Senate composition tally) down a row;
# fix up the SUM range so it again covers the tally rows above it.
$ws.Range("N39").Formula = "=SUM(N32:N38)"

# --- Minor sheet formatting / view tweaks ---

# Column D got a bit wider.
$ws.Columns.Item(4).ColumnWidth = 11.15

# Leave the selection on the newly added row.
$ws.Range("E44").Select()
